# Update cryptocurrency price (D) and volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.Formula = "'" + $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "44.182.51"
Set-TextCell $ws.Range("E2") "  +1.25%  "
Set-TextCell $ws.Range("D3") "2.245.31"
Set-TextCell $ws.Range("E3") "  +0.89%  "
Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  +0.15%  "
Set-TextCell $ws.Range("D5") "306.50"
Set-TextCell $ws.Range("E5") "  -2.27%  "
Set-TextCell $ws.Range("D6") "95.54"
Set-TextCell $ws.Range("E6") "  -2.69%  "
Set-TextCell $ws.Range("D7") "0.573"
Set-TextCell $ws.Range("E7") "  +0.71%  "
Set-TextCell $ws.Range("E8") "  +0.20%  "
Set-TextCell $ws.Range("E9") "  -1.56%  "
Set-TextCell $ws.Range("D10") "35.11"
Set-TextCell $ws.Range("E10") "  -2.68%  "
Set-TextCell $ws.Range("D11") "0.0814"
Set-TextCell $ws.Range("E11") "  -0.94%  "
Set-TextCell $ws.Range("D12") "7.22"
Set-TextCell $ws.Range("E12") "  -2.10%  "
Set-TextCell $ws.Range("E13") "  +0.07%  "
Set-TextCell $ws.Range("D14") "2.587.46"
Set-TextCell $ws.Range("E14") "  +0.94%  "
Set-TextCell $ws.Range("D15") "2.329.36"
Set-TextCell $ws.Range("E15") "  +4.77%  "
Set-TextCell $ws.Range("D16") "0.833"
Set-TextCell $ws.Range("E16") "  -0.72%  "
Set-TextCell $ws.Range("D17") "13.61"
Set-TextCell $ws.Range("E17") "  -3.36%  "
Set-TextCell $ws.Range("D18") "44.038.24"
Set-TextCell $ws.Range("E18") "  +1.21%  "
Set-TextCell $ws.Range("D19") "0.0₃0969"
Set-TextCell $ws.Range("E19") "  +0.49%  "
Set-TextCell $ws.Range("E20") "  +1.42%  "
Set-TextCell $ws.Range("D21") "12.12"
Set-TextCell $ws.Range("E21") "  -6.94%  "
Set-TextCell $ws.Range("D22") "65.44"
Set-TextCell $ws.Range("E22") "  +0.08%  "
Set-TextCell $ws.Range("D23") "236.56"
Set-TextCell $ws.Range("E23") "  +0.57%  "
Set-TextCell $ws.Range("D24") "2.96"
Set-TextCell $ws.Range("E24") "  -1.04%  "
Set-TextCell $ws.Range("D25") "2.00"
Set-TextCell $ws.Range("E25") "  -1.63%  "
Set-TextCell $ws.Range("E26") "  -0.02%  "
Set-TextCell $ws.Range("D27") "9.98"
Set-TextCell $ws.Range("E27") "  -0.38%  "
Set-TextCell $ws.Range("E28") "  -0.87%  "
Set-TextCell $ws.Range("D29") "37.49"
Set-TextCell $ws.Range("E29") "  +3.35%  "
Set-TextCell $ws.Range("D30") "5.99"
Set-TextCell $ws.Range("E30") "  +0.61%  "
Set-TextCell $ws.Range("D31") "20.15"
Set-TextCell $ws.Range("E31") "  +1.30%  "
Set-TextCell $ws.Range("D32") "153.34"
Set-TextCell $ws.Range("E32") "  -3.81%  "
Set-TextCell $ws.Range("D33") "0.0801"
Set-TextCell $ws.Range("E33") "  -3.39%  "
Set-TextCell $ws.Range("D34") "3.26"
Set-TextCell $ws.Range("E34") "  +2.90%  "
Set-TextCell $ws.Range("D35") "2.60"
Set-TextCell $ws.Range("E35") "  -2.93%  "
Set-TextCell $ws.Range("E36") "  +3.07%  "
Set-TextCell $ws.Range("E37") "  -0.99%  "
Set-TextCell $ws.Range("D38") "1.75"
Set-TextCell $ws.Range("E38") "  -7.10%  "
Set-TextCell $ws.Range("D39") "3.46"
Set-TextCell $ws.Range("E39") "  -3.12%  "
Set-TextCell $ws.Range("D40") "3.87"
Set-TextCell $ws.Range("E40") "  -3.15%  "
Set-TextCell $ws.Range("D41") "14.60"
Set-TextCell $ws.Range("E41") "  -7.03%  "
Set-TextCell $ws.Range("E42") "  -2.96%  "
Set-TextCell $ws.Range("E43") "  +0.26%  "
Set-TextCell $ws.Range("D44") "1.740.58"
Set-TextCell $ws.Range("E44") "  +1.90%  "
Set-TextCell $ws.Range("D45") "83.06"
Set-TextCell $ws.Range("E45") "  +1.20%  "
Set-TextCell $ws.Range("E46") "  -1.92%  "
Set-TextCell $ws.Range("E47") "  +6.52%  "
Set-TextCell $ws.Range("D48") "100.08"
Set-TextCell $ws.Range("E48") "  -1.51%  "
Set-TextCell $ws.Range("E49") "  -2.84%  "
Set-TextCell $ws.Range("D50") "8.17"
Set-TextCell $ws.Range("E50") "  +1.76%  "
Set-TextCell $ws.Range("E51") "  -3.02%  "
